$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-11-17 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-11-18 Tuesday", 2) | Out-Null

# Update each arithmetic-problem cell in the practice table
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "43+13="  # was "45+41="
$t.Cell(1, 2).Range.Text = "54+5="  # was "75+21="
$t.Cell(1, 3).Range.Text = "72-15="  # was "95-57="
$t.Cell(1, 4).Range.Text = "15-6="  # was "56-29="
$t.Cell(1, 5).Range.Text = "58-57="  # was "58+28="
$t.Cell(2, 1).Range.Text = "60-55="  # was "11+27="
$t.Cell(2, 2).Range.Text = "53+42="  # was "34+7="
$t.Cell(2, 3).Range.Text = "84-74="  # was "51+11="
$t.Cell(2, 4).Range.Text = "67-10="  # was "90-31="
$t.Cell(2, 5).Range.Text = "33+40="  # was "26+39="
$t.Cell(3, 1).Range.Text = "11+3="  # was "34+47="
$t.Cell(3, 2).Range.Text = "47+25="  # was "0+52="
$t.Cell(3, 3).Range.Text = "90-56="  # was "40+17="
$t.Cell(3, 4).Range.Text = "68+16="  # was "21+38="
$t.Cell(3, 5).Range.Text = "40+35="  # was "73-33="
$t.Cell(4, 1).Range.Text = "34+65="  # was "71-52="
$t.Cell(4, 2).Range.Text = "94-46="  # was "85-74="
$t.Cell(4, 3).Range.Text = "63+4="  # was "14+79="
$t.Cell(4, 4).Range.Text = "33+30="  # was "58+38="
$t.Cell(4, 5).Range.Text = "53+5="  # was "65-56="
$t.Cell(5, 1).Range.Text = "21+28="  # was "72-36="
$t.Cell(5, 2).Range.Text = "26+65="  # was "40+12="
$t.Cell(5, 3).Range.Text = "75-0="  # was "45+52="
$t.Cell(5, 4).Range.Text = "28+31="  # was "8+9="
$t.Cell(5, 5).Range.Text = "85+7="  # was "44-16="
$t.Cell(6, 1).Range.Text = "79-43="  # was "72-60="
$t.Cell(6, 2).Range.Text = "58-17="  # was "42+8="
$t.Cell(6, 3).Range.Text = "36+59="  # was "54-13="
$t.Cell(6, 4).Range.Text = "10+21="  # was "88-77="
$t.Cell(6, 5).Range.Text = "24-9="  # was "53+26="
$t.Cell(7, 1).Range.Text = "34+50="  # was "98-14="
$t.Cell(7, 2).Range.Text = "94-53="  # was "55+41="
$t.Cell(7, 3).Range.Text = "81-31="  # was "38+3="
$t.Cell(7, 4).Range.Text = "63+28="  # was "96-9="
$t.Cell(7, 5).Range.Text = "35+27="  # was "12+74="
$t.Cell(8, 1).Range.Text = "90-35="  # was "54+4="
$t.Cell(8, 2).Range.Text = "54+34="  # was "60-52="
$t.Cell(8, 3).Range.Text = "60-57="  # was "14+45="
$t.Cell(8, 4).Range.Text = "11+69="  # was "2+55="
$t.Cell(8, 5).Range.Text = "12+7="  # was "43+46="
$t.Cell(9, 1).Range.Text = "7+68="  # was "64+25="
$t.Cell(9, 2).Range.Text = "51-16="  # was "59-17="
$t.Cell(9, 3).Range.Text = "67-27="  # was "53-28="
$t.Cell(9, 4).Range.Text = "32-19="  # was "48-33="
$t.Cell(9, 5).Range.Text = "43+33="  # was "52-31="
$t.Cell(10, 1).Range.Text = "94-72="  # was "64+16="
$t.Cell(10, 2).Range.Text = "84+2="  # was "78-21="
$t.Cell(10, 3).Range.Text = "89-0="  # was "39-4="
$t.Cell(10, 4).Range.Text = "95-1="  # was "16+12="
$t.Cell(10, 5).Range.Text = "34+14="  # was "36+57="
$t.Cell(11, 1).Range.Text = "20+66="  # was "77-32="
$t.Cell(11, 2).Range.Text = "32+42="  # was "42+52="
$t.Cell(11, 3).Range.Text = "38-32="  # was "83-52="
$t.Cell(11, 4).Range.Text = "93-63="  # was "90+6="
$t.Cell(11, 5).Range.Text = "70+29="  # was "93-25="
$t.Cell(12, 1).Range.Text = "53+44="  # was "19+31="
$t.Cell(12, 2).Range.Text = "25+38="  # was "59+7="
$t.Cell(12, 3).Range.Text = "41-31="  # was "82-26="
$t.Cell(12, 4).Range.Text = "50+27="  # was "35+49="
$t.Cell(12, 5).Range.Text = "97-28="  # was "6+12="
$t.Cell(13, 1).Range.Text = "3+75="  # was "26+8="
$t.Cell(13, 2).Range.Text = "36-0="  # was "1+90="
$t.Cell(13, 3).Range.Text = "11+26="  # was "10+48="
$t.Cell(13, 4).Range.Text = "43+12="  # was "9+38="
$t.Cell(13, 5).Range.Text = "83+16="  # was "60-46="
$t.Cell(14, 1).Range.Text = "96-35="  # was "56+18="
$t.Cell(14, 2).Range.Text = "12+63="  # was "76-40="
$t.Cell(14, 3).Range.Text = "45-29="  # was "56+2="
$t.Cell(14, 4).Range.Text = "82-52="  # was "25-2="
$t.Cell(14, 5).Range.Text = "58-2="  # was "50-29="
$t.Cell(15, 1).Range.Text = "38+17="  # was "24+56="
$t.Cell(15, 2).Range.Text = "70-16="  # was "8+14="
$t.Cell(15, 3).Range.Text = "30+15="  # was "70+10="
$t.Cell(15, 4).Range.Text = "63-11="  # was "31+60="
$t.Cell(15, 5).Range.Text = "6-3="  # was "55-20="
$t.Cell(16, 1).Range.Text = "55-41="  # was "76-23="
$t.Cell(16, 2).Range.Text = "14+11="  # was "19-1="
$t.Cell(16, 3).Range.Text = "61-34="  # was "92-87="
$t.Cell(16, 4).Range.Text = "81-41="  # was "13+72="
$t.Cell(16, 5).Range.Text = "6+88="  # was "67-32="
$t.Cell(17, 1).Range.Text = "84-25="  # was "80-10="
$t.Cell(17, 2).Range.Text = "88-38="  # was "31-4="
$t.Cell(17, 3).Range.Text = "60-16="  # was "38+61="
$t.Cell(17, 4).Range.Text = "7+41="  # was "56+25="
$t.Cell(17, 5).Range.Text = "31-8="  # was "84-80="
$t.Cell(18, 1).Range.Text = "80-35="  # was "13+79="
$t.Cell(18, 2).Range.Text = "45-7="  # was "90-1="
$t.Cell(18, 3).Range.Text = "37-8="  # was "39+27="
$t.Cell(18, 4).Range.Text = "27+51="  # was "54+2="
$t.Cell(18, 5).Range.Text = "80-13="  # was "97-38="
$t.Cell(19, 1).Range.Text = "74-63="  # was "20+72="
$t.Cell(19, 2).Range.Text = "59-45="  # was "90-48="
$t.Cell(19, 3).Range.Text = "45+11="  # was "84-27="
$t.Cell(19, 4).Range.Text = "78-25="  # was "70-44="
$t.Cell(19, 5).Range.Text = "73+20="  # was "41+38="
$t.Cell(20, 1).Range.Text = "1+88="  # was "97-38="
$t.Cell(20, 2).Range.Text = "81-13="  # was "6+83="
$t.Cell(20, 3).Range.Text = "11+84="  # was "56+33="
$t.Cell(20, 4).Range.Text = "70-41="  # was "26+4="
$t.Cell(20, 5).Range.Text = "3+75="  # was "13+44="
